# Sprint 1 Backlog - "Minor Updates to remaining hours"
# Fill in the daily (Mon-Fri) remaining-hours tracker for the two task
# tables on Sheet1: add the missing Wed/Thurs/Fri (G:I) values for every
# task row, and backfill the missing Mon (E) value for the second table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Table 1: "1. Debug existing application" (rows 6-10) -----------------
# Mon (E) and Tues (F) already had values; add Wed/Thurs/Fri = 0.
$ws.Range("G6:I10").Value = 0

# --- Table 2: "2. Sorting Shares" (rows 13-17) ------------------------------
# Tues (F) already had a value; backfill Mon (E) = 1 and add
# Wed/Thurs/Fri = 0 for each task row.
$ws.Range("E13:E17").Value = 1
$ws.Range("G13:I17").Value = 0

# Update the last used cell/selection to reflect where data entry ended.
$ws.Range("E18").Select()
